# Apply the edits described by the diff to the active document.
$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Greeting: "Dear Mr. Liu," -> "Dear Sir / Dear Madam,"
#    The original paragraph is three runs: "Dear Mr. " | "Liu" | ",".
#    The target keeps three runs too, just with different wording in the
#    first two: "Dear " | "Sir / Dear Madam" | ",". A plain Find/Replace
#    would coalesce the (now identically-formatted) runs into one, so we
#    round-trip the paragraph's OOXML via WordOpenXML/InsertXML, editing the
#    text nodes directly while leaving the run/element structure intact.
# ---------------------------------------------------------------------------
$greeting = $d.Paragraphs.Item(1)
$gr = $greeting.Range
if ($gr.Text -like "Dear Mr. Liu*") {
    $gxml = $gr.WordOpenXML
    $gxml = $gxml.Replace('>Dear Mr. </w:t>', '>Dear </w:t>')
    $gxml = $gxml.Replace('>Liu</w:t>', '>Sir / Dear Madam</w:t>')
    $gr.InsertXML($gxml)
}

# ---------------------------------------------------------------------------
# 2) "Analytical skills gained through Bachelor's degree in physics"
#    The three runs that used to surround the grammar-check markers merge
#    into a single run (and the <w:proofErr/> pair is dropped).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Analytical skills gained through Bachelor" + [char]0x2019 + "s degree in physics", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Analytical skills gained through Bachelor" + [char]0x2019 + "s degree in physics", 2)

# ---------------------------------------------------------------------------
# 3) "Ability to develop HTML templates from mockups that are both
#    pixel-perfect and responsive" - same kind of merge as #2, removing the
#    spell-check markers around "mockups".
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Ability to develop HTML templates from mockups that are both pixel-perfect and responsive", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Ability to develop HTML templates from mockups that are both pixel-perfect and responsive", 2)

# ---------------------------------------------------------------------------
# 4) "Hyungmo Gu" - the two runs that used to surround the spell-check
#    markers merge into a single run.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Hyungmo Gu", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Hyungmo Gu", 2)

# The Find/Replace above merges the runs and drops the proofErr markers that
# fall inside the matched range, but a leading <w:proofErr w:type="spellStart"/>
# that sits right at the very start of the "Hyungmo Gu" paragraph (before any
# matched text) is not part of the replaced range, so it survives. Round-trip
# that paragraph's XML through WordOpenXML/InsertXML to normalize it away
# while preserving the run formatting.
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*Hyungmo Gu*") {
        $r = $p.Range
        $xml = $r.WordOpenXML
        $r.InsertXML($xml)
        break
    }
}
